$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8104699850082397
$ws.Range("B1").Value = 1.22692883014679
$ws.Range("C1").Value = 2.458017349243164
$ws.Range("D1").Value = 3.739233016967773
$ws.Range("E1").Value = 2.731669664382935
